# Apply updated cryptocurrency market data to Sheet1 (coin name/link swaps + price/volume refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    # Force text storage so purely-numeric-looking strings (e.g. "1.002") are not
    # auto-coerced into numbers by Excel, matching the source data's inline-string type.
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
}

$ws.Range("D2").Value2 = "29.373.90"
$ws.Range("E2").Value2 = "  +0.09%  "
$ws.Range("D3").Value2 = "1.881.95"
$ws.Range("E3").Value2 = "  +0.27%  "
Set-TextCell "D4" "1.002"
$ws.Range("E4").Value2 = "  +0.16%  "
Set-TextCell "D5" "0.7105"
$ws.Range("E5").Value2 = "  -0.18%  "
Set-TextCell "D6" "242.90"
$ws.Range("E6").Value2 = "  +0.21%  "
$ws.Range("E7").Value2 = "  +0.11%  "
Set-TextCell "D8" "0.08013"
$ws.Range("E8").Value2 = "  +3.42%  "
$ws.Range("E9").Value2 = "  +0.77%  "
Set-TextCell "D10" "25.14"
$ws.Range("E10").Value2 = "  +0.81%  "
Set-TextCell "D11" "0.08348"
$ws.Range("E11").Value2 = "  -2.12%  "
$ws.Range("D12").Value2 = "1.890.79"
$ws.Range("E12").Value2 = "  +0.62%  "
Set-TextCell "D13" "5.246"
$ws.Range("E13").Value2 = "  +0.33%  "
$ws.Range("B14").Value2 = "Polygon"
$ws.Range("C14").Value2 = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell "D14" "0.7181"
$ws.Range("E14").Value2 = "  +1.03%  "
$ws.Range("B15").Value2 = "Litecoin"
$ws.Range("C15").Value2 = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D15" "94.16"
$ws.Range("E15").Value2 = "  +2.76%  "
Set-TextCell "D16" "6.329"
$ws.Range("E16").Value2 = "  +5.05%  "
Set-TextCell "D17" "0.000008527"
$ws.Range("E17").Value2 = "  +3.54%  "
$ws.Range("D18").Value2 = "29.394.74"
Set-TextCell "D19" "241.94"
$ws.Range("E19").Value2 = "  +0.30%  "
$ws.Range("D20").Value2 = "2.132.26"
$ws.Range("E20").Value2 = "  +0.30%  "
Set-TextCell "D21" "13.24"
$ws.Range("E21").Value2 = "  -0.22%  "
$ws.Range("E22").Value2 = "  +0.08%  "
Set-TextCell "D23" "7.821"
$ws.Range("E23").Value2 = "  +0.43%  "
Set-TextCell "D24" "1.002"
$ws.Range("E24").Value2 = "  +0.19%  "
Set-TextCell "D25" "0.1583"
$ws.Range("E25").Value2 = "  -0.66%  "
Set-TextCell "D26" "163.43"
$ws.Range("E26").Value2 = "  +0.03%  "
Set-TextCell "D27" "9.084"
$ws.Range("E27").Value2 = "  +0.47%  "
Set-TextCell "D28" "18.61"
$ws.Range("E28").Value2 = "  +0.68%  "
$ws.Range("E29").Value2 = "  -0.20%  "
Set-TextCell "D30" "4.412"
$ws.Range("E30").Value2 = "  +0.22%  "
Set-TextCell "D31" "4.328"
$ws.Range("E31").Value2 = "  -1.10%  "
Set-TextCell "D32" "1.201"
Set-TextCell "D33" "0.05372"
$ws.Range("E33").Value2 = "  +1.75%  "
Set-TextCell "D34" "1.944"
$ws.Range("E34").Value2 = "  +0.25%  "
$ws.Range("B35").Value2 = "ImmutableX"
$ws.Range("C35").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell "D35" "0.7749"
$ws.Range("E35").Value2 = "  +3.18%  "
$ws.Range("B36").Value2 = "ARBITRUM"
$ws.Range("C36").Value2 = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell "D36" "1.181"
$ws.Range("E36").Value2 = "  +0.34%  "
Set-TextCell "D37" "2.697"
$ws.Range("E37").Value2 = "  +0.35%  "
Set-TextCell "D38" "0.01884"
$ws.Range("E38").Value2 = "  +0.83%  "
$ws.Range("D39").Value2 = "1.282.38"
$ws.Range("E39").Value2 = "  +8.12%  "
Set-TextCell "D40" "2.749"
$ws.Range("E40").Value2 = "  +1.00%  "
Set-TextCell "D41" "6.552"
$ws.Range("E41").Value2 = "  +2.01%  "
Set-TextCell "D42" "0.9190"
$ws.Range("E42").Value2 = "  +3.99%  "
$ws.Range("B43").Value2 = "Quant"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell "D43" "112.95"
$ws.Range("E43").Value2 = "  +5.49%  "
$ws.Range("B44").Value2 = "Aave"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D44" "74.56"
$ws.Range("E44").Value2 = "  +2.14%  "
$ws.Range("E45").Value2 = "  +0.10%  "
Set-TextCell "D46" "0.00000000129"
$ws.Range("E46").Value2 = "  +5.69%  "
$ws.Range("D47").Value2 = "2.024.72"
$ws.Range("E47").Value2 = "  -0.15%  "
$ws.Range("E48").Value2 = "  -0.45%  "
Set-TextCell "D49" "0.5225"
$ws.Range("E49").Value2 = "  +0.30%  "
Set-TextCell "D50" "9.536"
$ws.Range("E50").Value2 = "  +1.18%  "
Set-TextCell "D51" "0.4376"
$ws.Range("E51").Value2 = "  +1.30%  "
